# "Yearly coverage in scenario 1"
#
# The "Platform Coverage" sheet previously only had a coverage value of 0.6
# every other year (2018, 2020, 2022, ...) for the first MDA/School row.
# This fills in the in-between (odd) years so every year from 2018-2040
# has a coverage value of 0.6.

$wb = $excel.ActiveWorkbook
$wsCoverage = $wb.Worksheets.Item("Platform Coverage")
$wsMarket   = $wb.Worksheets.Item("MarketShare")

# Columns H:AD on row 2 correspond to years 2018-2040. Only the even years
# (H, J, L, N, P, R, T, V, X, Z, AB, AD) were populated; fill in the
# remaining odd years with the same 0.6 coverage value.
$yearCols = @("I","K","M","O","Q","S","U","W","Y","AA","AC")
foreach ($col in $yearCols) {
    $wsCoverage.Range($col + "2").Value = 0.6
}

# Restore the selections/active sheet seen when the file was last saved.
$wsCoverage.Range("AE2").Select()
$wsMarket.Activate()
$wsMarket.Range("Z3").Select()
